$wb = $excel.ActiveWorkbook

# --- Fill in row 12 on the PLA sheet with the new "Yellow" filament preset ---
$plaSheet = $wb.Worksheets.Item("PLA")

$plaSheet.Range("A12").Value = 11
$plaSheet.Range("B12").Value = "3DPlast"
$plaSheet.Range("C12").Value = "PLA"
$plaSheet.Range("D12").Value = "Yellow"
$plaSheet.Range("E12").Value = 190
$plaSheet.Range("F12").Value = 220
$plaSheet.Range("G12").Value = 220
$plaSheet.Range("H12").Value = 220
$plaSheet.Range("I12").Value = 0.963
$plaSheet.Range("J12").Value = 0.032

# --- Update the active selection on each sheet ---
$plaSheet.Range("N11").Select() | Out-Null

$petgSheet = $wb.Worksheets.Item("PETG")
$petgSheet.Range("O10").Select() | Out-Null

# --- Make PLA the active/selected tab (tabSelected moves from PETG to PLA) ---
$excel.Worksheets.Item("PLA").Activate() | Out-Null
